$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: C3 changes from "test-id-14" to "test-id-15" (A3/B3 unchanged)
$ws.Range("C3").Value = "test-id-15"

# Rows 4-6: clear the CODE/DATE/ADDER values entirely (B keeps its date
# number-format style but becomes blank; A and C become fully empty)
$ws.Range("A4:C6").ClearContents()

# Update the active selection to just cell C4
[void]$ws.Range("C4").Select()
